$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 24103.8
$ws.Range("I21").Value = 23833.334
$ws.Range("J21").Value = 24509.5
$ws.Range("K21").Value = 23833.334
$ws.Range("L21").Value = 24509.5
$ws.Range("M21").Value = -23365.334
$ws.Range("N21").Value = -25445.5

$ws.Range("H23").Value = 24103.8
$ws.Range("I23").Value = 23833.334
$ws.Range("J23").Value = 24509.5
$ws.Range("K23").Value = 23833.334
$ws.Range("L23").Value = 24509.5
$ws.Range("M23").Value = -23599.334
$ws.Range("N23").Value = -24977.5

$ws.Range("H38").Value = 1471.5454
$ws.Range("J38").Value = 2863.6365
$ws.Range("L38").Value = 8590.9095
$ws.Range("N38").Value = -9334.9095

$ws.Range("H94").Value = 5063.4165
$ws.Range("I94").Value = 5063.4165
$ws.Range("K94").Value = 5063.4165
$ws.Range("M94").Value = -4612.4165

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 25665
$ws.Range("I21").Value = 25665
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 25665
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -25291
$ws.Range("N21").ClearContents()

$ws.Range("H44").Value = 39900
$ws.Range("J44").Value = 39900
$ws.Range("L44").Value = 39900
$ws.Range("N44").Value = -40876

$ws.Range("H105").Value = 500370
$ws.Range("J105").Value = 500370
$ws.Range("L105").Value = 500370
$ws.Range("N105").Value = -507358

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 37750
$ws.Range("J95").Value = 37750
$ws.Range("L95").Value = 37750
$ws.Range("N95").Value = -43242

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1660.2
$ws.Range("I16").Value = 1575
$ws.Range("J16").Value = 2001
$ws.Range("K16").Value = 1575
$ws.Range("L16").Value = 2001
$ws.Range("M16").Value = -1288
$ws.Range("N16").Value = -2575

$ws.Range("H32").Value = 5500
$ws.Range("J32").Value = 10000
$ws.Range("L32").Value = 10000
$ws.Range("N32").Value = -10632

$ws.Range("H96").Value = 19994.8
$ws.Range("J96").Value = 19994.8
$ws.Range("L96").Value = 19994.8
$ws.Range("N96").Value = -25486.8

$ws.Range("H113").Value = 1660.2
$ws.Range("I113").Value = 1575
$ws.Range("J113").Value = 2001
$ws.Range("K113").Value = 1575
$ws.Range("L113").Value = 2001
$ws.Range("M113").Value = 595
$ws.Range("N113").Value = -6341

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 235.71428
$ws.Range("I26").Value = 191.66667
$ws.Range("J26").Value = 500
$ws.Range("K26").Value = 575.00001
$ws.Range("L26").Value = 1500
$ws.Range("M26").Value = -287.00001
$ws.Range("N26").Value = -2076

$ws.Range("H51").Value = 1600
$ws.Range("I51").Value = 1000
$ws.Range("J51").Value = 1800
$ws.Range("K51").Value = 3000
$ws.Range("L51").Value = 5400
$ws.Range("M51").Value = -2540
$ws.Range("N51").Value = -6320

$ws.Range("H98").Value = 167607.83
$ws.Range("I98").Value = 250360.75
$ws.Range("J98").Value = 2102
$ws.Range("K98").Value = 751082.25
$ws.Range("L98").Value = 6306
$ws.Range("M98").Value = -749584.25
$ws.Range("N98").Value = -9302

$ws.Range("H99").Value = 5150

$ws.Range("H100").Value = 50000
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws.Range("H101").Value = 6147.5454
$ws.Range("J101").Value = 6147.5454
$ws.Range("L101").Value = 18442.6362
$ws.Range("N101").Value = -23310.6362

$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").ClearContents()

$ws.Range("H103").Value = 2564.6
$ws.Range("I103").Value = 641.3333
$ws.Range("J103").Value = 5449.5
$ws.Range("K103").Value = 1923.9999
$ws.Range("L103").Value = 16348.5
$ws.Range("M103").Value = -1044.9999
$ws.Range("N103").Value = -18106.5

$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws.Range("H110").Value = 13333.333
$ws.Range("I110").Value = 10000
$ws.Range("K110").Value = 30000
$ws.Range("M110").Value = -25910

$ws.Range("H111").Value = 560
$ws.Range("I111").Value = 120
$ws.Range("J111").Value = 1000
$ws.Range("K111").Value = 360
$ws.Range("L111").Value = 3000
$ws.Range("M111").Value = 2707
$ws.Range("N111").Value = -9134

$ws.Range("H112").Value = 4009
$ws.Range("J112").Value = 4009
$ws.Range("L112").Value = 12027
$ws.Range("N112").Value = -14243

$ws.Range("H113").Value = 714928.7
$ws.Range("J113").Value = 664.4545000000001
$ws.Range("L113").Value = 1993.3635
$ws.Range("N113").Value = -6333.3635

$ws.Range("H121").Value = 1209.5454
$ws.Range("I121").Value = 706.5
$ws.Range("J121").Value = 1813.2
$ws.Range("K121").Value = 2119.5
$ws.Range("L121").Value = 5439.6
$ws.Range("M121").Value = -809.5
$ws.Range("N121").Value = -8059.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()

$ws.Range("H98").Value = 28000
$ws.Range("J98").Value = 28000
$ws.Range("L98").Value = 28000
$ws.Range("N98").Value = -33990

$ws.Range("H109").Value = 11823.308
$ws.Range("J109").Value = 11823.308
$ws.Range("L109").Value = 11823.308
$ws.Range("N109").Value = -13903.308

$ws.Range("H113").Value = 1508.381
$ws.Range("I113").Value = 945.1
$ws.Range("J113").Value = 2020.4546
$ws.Range("K113").Value = 945.1
$ws.Range("L113").Value = 2020.4546
$ws.Range("M113").Value = 1224.9
$ws.Range("N113").Value = -6360.4546

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 42475
$ws.Range("I34").Value = 42475
$ws.Range("K34").Value = 42475
$ws.Range("M34").Value = -42303

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 25150
$ws.Range("I8").Value = 500
$ws.Range("K8").Value = 500
$ws.Range("M8").Value = -360

$ws.Range("H23").Value = 2499
$ws.Range("I23").Value = 2499
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 2499
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -2270
$ws.Range("N23").ClearContents()

$ws.Range("H100").Value = 19932.562
$ws.Range("I100").Value = 34101.668
$ws.Range("J100").Value = 1715.1428
$ws.Range("K100").Value = 68203.336
$ws.Range("L100").Value = 3430.2856
$ws.Range("M100").Value = -4512.2856
